# "Generate Report for Handoff" - refresh the handoff/handback report data.
#
# The "Handed back: in sync with en-US" status is stale; the files are now
# ready to be handed off again, and the latest-handback freshness check has
# flagged that the checked-in handback files are behind the latest commit.
# Update the Overview + per-locale sheets accordingly.

$wb = $excel.ActiveWorkbook

$readyForHandoff = "Ready for handoff"
$overviewDate    = "2016-09-06 03:22:42"
$zhcnHandoffDate = "2016-09-06 03:22:31"

$notLatest45f0 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b0006ec4aed81ba9a7579ac358b423baf39ab6b/e2e/45f0adaa-5bc7-45b4-a7af-1aceb5614af9.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cac5a39a03d82d8d9cb5ab04b6eca44b1b4f7e18/e2e/45f0adaa-5bc7-45b4-a7af-1aceb5614af9.md."
$notLatest7eb0 = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/3b0006ec4aed81ba9a7579ac358b423baf39ab6b/e2e/7eb0bc2f-3120-4fa9-a9cf-9f7c8875b9e0.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/cac5a39a03d82d8d9cb5ab04b6eca44b1b4f7e18/e2e/7eb0bc2f-3120-4fa9-a9cf-9f7c8875b9e0.md."

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Range("E2").Value = $readyForHandoff
$wsOverview.Range("F2").Value = $readyForHandoff
$wsOverview.Range("G2").Value = $overviewDate

$wsOverview.Range("E3").Value = $readyForHandoff
$wsOverview.Range("F3").Value = $readyForHandoff
$wsOverview.Range("G3").Value = $overviewDate

$wsOverview.Columns.Item(5).ColumnWidth = 16.33
$wsOverview.Columns.Item(6).ColumnWidth = 16.33

# ---- zh-cn sheet ----
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Range("C2").Value = $readyForHandoff
$wsZhCn.Range("H2").Value = $zhcnHandoffDate
$wsZhCn.Range("P2").Value = $notLatest45f0

$wsZhCn.Range("C3").Value = $readyForHandoff
$wsZhCn.Range("H3").Value = $zhcnHandoffDate
$wsZhCn.Range("P3").Value = $notLatest7eb0

$wsZhCn.Columns.Item(3).ColumnWidth = 16.33
$wsZhCn.Columns.Item(16).ColumnWidth = 39.17

# ---- de-de sheet ----
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Range("C2").Value = $readyForHandoff
$wsDeDe.Range("H2").Value = $overviewDate
$wsDeDe.Range("P2").Value = $notLatest45f0

$wsDeDe.Range("C3").Value = $readyForHandoff
$wsDeDe.Range("H3").Value = $overviewDate
$wsDeDe.Range("P3").Value = $notLatest7eb0

$wsDeDe.Columns.Item(3).ColumnWidth = 16.33
$wsDeDe.Columns.Item(16).ColumnWidth = 39.17
